$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet1: rewrite the condition table. Column D ("opacity"/0) is gone,
# the remaining A:C columns get new values ("images/..." paths and a
# single merged question string instead of "Question?"/"Rest").
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$data = @(
    @("arrow", "image",            "q_or_r"),
    @("down",  "images/down.jpg",  "Rate your ability to control your brain"),
    @("up",    "images/up.jpg",    "Rate your ability to control your brain"),
    @("up",    "images/up.jpg",    "Rate your ability to control your brain"),
    @("down",  "images/down.jpg",  "Rate your ability to control your brain"),
    @("down",  "images/down.jpg",  "Rate your ability to control your brain"),
    @("up",    "images/up.jpg",    "Rate your ability to control your brain")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Drop the old column D entirely (now unused).
[void]$ws1.Columns.Item(4).Delete()

# Column B now holds longer strings ("images/....jpg") so widen it;
# columns A and C keep the sheet's normal width.
$ws1.Columns.Item(2).ColumnWidth = 14.3

# ---------------------------------------------------------------------
# Sheet2 / Sheet3: unchanged data, only the saved selection differs
# (same B7:C11 block carried over from Sheet1 plus the sheet's own A1).
# ---------------------------------------------------------------------
foreach ($idx in 2, 3) {
    $ws = $wb.Worksheets.Item($idx)
    [void]$ws.Select()
    [void]$ws.Range("B7:C11").Select()
}

# Sheet1 stays the active/selected tab, with the new B7:C11 selection
# (the engine anchors the active cell at the selection's top-left
# corner).
[void]$ws1.Select()
[void]$ws1.Range("B7:C11").Select()
